# Apply the updated cryptocurrency snapshot values / symbol re-ordering
# as captured by the Tue Dec 27 05:16:32 UTC 2022 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$val0 = "'243.43"
$ws.Range("D2").Value = $val0

# Row 5
$val1 = "'0.05995"
$ws.Range("D5").Value = $val1

# Row 6
$val2 = "'3.425"
$ws.Range("D6").Value = $val2

# Row 7
$val3 = "'6.483"
$ws.Range("D7").Value = $val3

# Row 8
$val4 = "'0.8084"
$ws.Range("D8").Value = $val4

# Row 9
$val5 = "'0.9244"
$ws.Range("D9").Value = $val5

# Row 10
$val6 = "'0.1424"
$ws.Range("D10").Value = $val6

# Row 11
$val7 = "'0.07415"
$ws.Range("D11").Value = $val7

# Row 12
$val8 = "'0.03280"
$ws.Range("D12").Value = $val8

# Row 13
$val9 = "'0.03069"
$ws.Range("D13").Value = $val9

# Row 14
$val10 = "'0.09352"
$ws.Range("D14").Value = $val10

# Row 15
$val11 = "'3.848"
$ws.Range("D15").Value = $val11

# Row 16
$val12 = "'0.001577"
$ws.Range("D16").Value = $val12

# Row 17
$val13 = "'0.04699"
$ws.Range("D17").Value = $val13

# Row 18
$val14 = "'0.0005918"
$ws.Range("D18").Value = $val14

# Row 19
$val15 = "'0.005867"
$ws.Range("D19").Value = $val15

# Row 20
$val16 = "'0.001272"
$ws.Range("D20").Value = $val16
$val17 = "19BitKanKANBestin24h"
$ws.Range("E20").Value = $val17

# Row 21
$val18 = "'0.004881"
$ws.Range("D21").Value = $val18

# Row 22
$val19 = "'0.00006798"
$ws.Range("D22").Value = $val19

# Row 23
$val20 = "'3.566"
$ws.Range("D23").Value = $val20

# Row 24
$val21 = "'2.179"
$ws.Range("D24").Value = $val21

# Row 27
$val22 = "'0.0002652"
$ws.Range("D27").Value = $val22

# Row 40
$val23 = "'0.03971"
$ws.Range("D40").Value = $val23

# Row 41
$val24 = "BKEXToken"
$ws.Range("B41").Value = $val24
$val25 = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("C41").Value = $val25
$val26 = "'0.1079"
$ws.Range("D41").Value = $val26
$val27 = "40BKEXTokenBKK"
$ws.Range("E41").Value = $val27

# Row 42
$val28 = "'0.002649"
$ws.Range("D42").Value = $val28
$val29 = "41CEJICEJI"
$ws.Range("E42").Value = $val29

# Row 43
$val30 = "KickToken"
$ws.Range("B43").Value = $val30
$val31 = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("C43").Value = $val31
$val32 = "'0.003069"
$ws.Range("D43").Value = $val32
$val33 = "42KickTokenKICKWorstin24h"
$ws.Range("E43").Value = $val33

# Row 44
$val34 = "'0.009186"
$ws.Range("D44").Value = $val34

# Row 45
$val35 = "'0.00005069"
$ws.Range("D45").Value = $val35

# Row 47
$val36 = "'0.6498"
$ws.Range("D47").Value = $val36
$val37 = "46CoinbaseStockTokenCOIN"
$ws.Range("E47").Value = $val37

# Row 48
$val38 = "'0.002453"
$ws.Range("D48").Value = $val38

# Row 49
$val39 = "'0.00002099"
$ws.Range("D49").Value = $val39

# Row 50
$val40 = "'0.0001999"
$ws.Range("D50").Value = $val40
